$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 400.1111
$ws.Range("I55").Value = 366.83334
$ws.Range("J55").Value = 466.66666
$ws.Range("K55").Value = 366.83334
$ws.Range("L55").Value = 466.66666
$ws.Range("M55").Value = -152.83334
$ws.Range("N55").Value = -894.66666
# Row 70
$ws.Range("H70").Value = 128235.875
$ws.Range("I70").Value = 1871.25
$ws.Range("J70").Value = 254600.5
$ws.Range("K70").Value = 5613.75
$ws.Range("L70").Value = 763801.5
$ws.Range("M70").Value = -5343.75
$ws.Range("N70").Value = -764341.5
# Row 73
$ws.Range("H73").Value = 128235.875
$ws.Range("I73").Value = 1871.25
$ws.Range("J73").Value = 254600.5
$ws.Range("K73").Value = 5613.75
$ws.Range("L73").Value = 763801.5
$ws.Range("M73").Value = -4677.75
$ws.Range("N73").Value = -765673.5
# Row 98
$ws.Range("H98").Value = 4317.4736
$ws.Range("I98").Value = 1403.3636
$ws.Range("J98").Value = 8324.375
$ws.Range("K98").Value = 1403.3636
$ws.Range("L98").Value = 8324.375
$ws.Range("M98").Value = 94.63640000000009
$ws.Range("N98").Value = -11320.375
# Row 106
$ws.Range("H106").Value = 8027.2856
$ws.Range("I106").Value = 8027.2856
$ws.Range("K106").Value = 8027.2856
$ws.Range("M106").Value = -7396.2856
# Row 113
$ws.Range("H113").Value = 1940
$ws.Range("I113").Value = 1940
$ws.Range("K113").Value = 1940
$ws.Range("M113").Value = 1314
# Row 116
$ws.Range("H116").Value = 5310.8667
$ws.Range("J116").Value = 9873
$ws.Range("L116").Value = 9873
$ws.Range("N116").Value = -16757
# Row 122
$ws.Range("H122").Value = 4317.4736
$ws.Range("I122").Value = 1403.3636
$ws.Range("J122").Value = 8324.375
$ws.Range("K122").Value = 4210.0908
$ws.Range("L122").Value = 24973.125
$ws.Range("M122").Value = -1760.0908
$ws.Range("N122").Value = -29873.125
# Row 132
$ws.Range("H132").Value = 2140.3286
$ws.Range("I132").Value = 2076.1343
$ws.Range("J132").Value = 3574
$ws.Range("K132").Value = 6228.402900000001
$ws.Range("L132").Value = 10722
$ws.Range("M132").Value = -3698.402900000001
$ws.Range("N132").Value = -15782
# Row 137
$ws.Range("H137").Value = 3594.6
$ws.Range("I137").Value = 3784.6667
$ws.Range("K137").Value = 11354.0001
$ws.Range("M137").Value = -8804.000100000001
# Row 138
$ws.Range("H138").Value = 3825.453
$ws.Range("I138").Value = 1802.2106
$ws.Range("J138").Value = 4956.0884
$ws.Range("K138").Value = 5406.6318
$ws.Range("L138").Value = 14868.2652
$ws.Range("M138").Value = -266.6318000000001
$ws.Range("N138").Value = -25148.2652
# Row 141
$ws.Range("H141").Value = 823.0909
$ws.Range("I141").Value = 823.0909
$ws.Range("K141").Value = 2469.2727
$ws.Range("M141").Value = 2710.7273

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9132.25
$ws.Range("I32").Value = 1049.0952
$ws.Range("J32").Value = 24563.727
$ws.Range("K32").Value = 1049.0952
$ws.Range("L32").Value = 24563.727
$ws.Range("M32").Value = -762.0952
$ws.Range("N32").Value = -25137.727
# Row 46
$ws.Range("H46").Value = 20497.2
$ws.Range("J46").Value = 20621.75
$ws.Range("L46").Value = 20621.75
$ws.Range("N46").Value = -21259.75
# Row 74
$ws.Range("H74").Value = 1699.875
$ws.Range("I74").Value = 1494.2941
$ws.Range("J74").Value = 2864.8333
$ws.Range("K74").Value = 1494.2941
$ws.Range("L74").Value = 2864.8333
$ws.Range("M74").Value = -620.2941000000001
$ws.Range("N74").Value = -4612.8333
# Row 77
$ws.Range("H77").Value = 1699.875
$ws.Range("I77").Value = 1494.2941
$ws.Range("J77").Value = 2864.8333
$ws.Range("K77").Value = 7471.4705
$ws.Range("L77").Value = 14324.1665
$ws.Range("M77").Value = -3103.4705
$ws.Range("N77").Value = -23060.1665
# Row 122
$ws.Range("H122").Value = 2929.0625
$ws.Range("I122").Value = 2984.4614
$ws.Range("J122").Value = 2689
$ws.Range("K122").Value = 8953.3842
$ws.Range("L122").Value = 8067
$ws.Range("M122").Value = -6503.3842
$ws.Range("N122").Value = -12967
# Row 125
$ws.Range("H125").Value = 89999
$ws.Range("J125").Value = 89999
$ws.Range("L125").Value = 89999
$ws.Range("N125").Value = -99839
# Row 132
$ws.Range("H132").Value = 2729.7173
$ws.Range("I132").Value = 2251.45
$ws.Range("K132").Value = 6754.349999999999
$ws.Range("M132").Value = -4224.349999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2482.8333
$ws.Range("I134").Value = 2430.25
$ws.Range("J134").Value = 2798.3333
$ws.Range("K134").Value = 7290.75
$ws.Range("L134").Value = 8394.999899999999
$ws.Range("M134").Value = -4755.75
$ws.Range("N134").Value = -13464.9999

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 15242.333
$ws.Range("I16").Value = 10211.177
$ws.Range("J16").Value = 36624.75
$ws.Range("K16").Value = 10211.177
$ws.Range("L16").Value = 36624.75
$ws.Range("M16").Value = -9924.177
$ws.Range("N16").Value = -37198.75
# Row 31
$ws.Range("H31").Value = 7009.3
$ws.Range("I31").Value = 6614.353
$ws.Range("J31").Value = 7525.769
$ws.Range("K31").Value = 6614.353
$ws.Range("L31").Value = 7525.769
$ws.Range("M31").Value = -6319.353
$ws.Range("N31").Value = -8115.769
# Row 34
$ws.Range("H34").Value = 7009.3
$ws.Range("I34").Value = 6614.353
$ws.Range("J34").Value = 7525.769
$ws.Range("K34").Value = 6614.353
$ws.Range("L34").Value = 7525.769
$ws.Range("M34").Value = -6412.353
$ws.Range("N34").Value = -7929.769
# Row 58
$ws.Range("H58").Value = 3828.375
$ws.Range("I58").Value = 4911.1816
$ws.Range("K58").Value = 4911.1816
$ws.Range("M58").Value = -4708.1816
# Row 113
$ws.Range("H113").Value = 15242.333
$ws.Range("I113").Value = 10211.177
$ws.Range("J113").Value = 36624.75
$ws.Range("K113").Value = 10211.177
$ws.Range("L113").Value = 36624.75
$ws.Range("M113").Value = -8041.177
$ws.Range("N113").Value = -40964.75
# Row 129
$ws.Range("H129").Value = 49999.668
$ws.Range("J129").Value = 49999.668
$ws.Range("L129").Value = 49999.668
$ws.Range("N129").Value = -59999.668
# Row 132
$ws.Range("H132").Value = 1305.65
$ws.Range("I132").Value = 1172.5098
$ws.Range("J132").Value = 2060.111
$ws.Range("K132").Value = 3517.5294
$ws.Range("L132").Value = 6180.333
$ws.Range("M132").Value = -987.5294000000004
$ws.Range("N132").Value = -11240.333
# Row 134
$ws.Range("H134").Value = 1363.1333
$ws.Range("I134").Value = 1350.7246
$ws.Range("J134").Value = 1505.8334
$ws.Range("K134").Value = 4052.1738
$ws.Range("L134").Value = 4517.5002
$ws.Range("M134").Value = -1517.1738
$ws.Range("N134").Value = -9587.5002
# Row 136
$ws.Range("H136").Value = 3828.375
$ws.Range("I136").Value = 4911.1816
$ws.Range("K136").Value = 14733.5448
$ws.Range("M136").Value = -12183.5448

$ws = $wb.Worksheets.Item("CUL")
# Row 128
$ws.Range("H128").Value = 156104
$ws.Range("I128").Value = 156104
$ws.Range("K128").Value = 468312
$ws.Range("M128").Value = -463332
# Row 129
$ws.Range("H129").Value = 1062.9231
$ws.Range("I129").Value = 616.9091
$ws.Range("J129").Value = 3516
$ws.Range("K129").Value = 1850.7273
$ws.Range("L129").Value = 10548
$ws.Range("M129").Value = 3149.2727
$ws.Range("N129").Value = -20548

$ws = $wb.Worksheets.Item("GSM")
# Row 63
$ws.Range("H63").Value = 17034
$ws.Range("I63").Value = 13103
$ws.Range("J63").Value = 18999.5
$ws.Range("K63").Value = 13103
$ws.Range("L63").Value = 18999.5
$ws.Range("M63").Value = -12417
$ws.Range("N63").Value = -20371.5
# Row 66
$ws.Range("H66").Value = 17034
$ws.Range("I66").Value = 13103
$ws.Range("J66").Value = 18999.5
$ws.Range("K66").Value = 39309
$ws.Range("L66").Value = 56998.5
$ws.Range("M66").Value = -35877
$ws.Range("N66").Value = -63862.5
# Row 122
$ws.Range("H122").Value = 4009.4
$ws.Range("I122").Value = 4062.111
$ws.Range("J122").Value = 3930.3333
$ws.Range("K122").Value = 12186.333
$ws.Range("L122").Value = 11790.9999
$ws.Range("M122").Value = -9736.332999999999
$ws.Range("N122").Value = -16690.9999
# Row 132
$ws.Range("H132").Value = 2433.054
$ws.Range("I132").Value = 2403.2917
$ws.Range("J132").Value = 3504.5
$ws.Range("K132").Value = 7209.875100000001
$ws.Range("L132").Value = 10513.5
$ws.Range("M132").Value = -4679.875100000001
$ws.Range("N132").Value = -15573.5
# Row 136
$ws.Range("H136").Value = 31019.572
$ws.Range("J136").Value = 31019.572
$ws.Range("L136").Value = 93058.716
$ws.Range("N136").Value = -98158.716

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1828.75
$ws.Range("I7").Value = 1828.75
$ws.Range("K7").Value = 1828.75
$ws.Range("M7").Value = -1716.75
# Row 120
$ws.Range("H120").Value = 8698
$ws.Range("J120").Value = 8698
$ws.Range("L120").Value = 8698
$ws.Range("N120").Value = -18374
# Row 126
$ws.Range("H126").Value = 1828.75
$ws.Range("I126").Value = 1828.75
$ws.Range("K126").Value = 5486.25
$ws.Range("M126").Value = -3016.25
# Row 130
$ws.Range("H130").Value = 119999.164
$ws.Range("J130").Value = 119999.164
$ws.Range("L130").Value = 119999.164
$ws.Range("N130").Value = -130039.164
# Row 132
$ws.Range("H132").Value = 6960
$ws.Range("I132").Value = 5633.087
$ws.Range("K132").Value = 16899.261
$ws.Range("M132").Value = -14369.261

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 3650
$ws.Range("I62").Value = 3650
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3650
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3026
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 3650
$ws.Range("I65").Value = 3650
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 18250
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -15130
$ws.Range("N65").ClearContents()
# Row 81
$ws.Range("H81").Value = 2590.8235
$ws.Range("I81").Value = 2185.7144
$ws.Range("J81").Value = 4481.3335
$ws.Range("K81").Value = 4371.4288
$ws.Range("L81").Value = 8962.666999999999
$ws.Range("M81").Value = -3310.4288
$ws.Range("N81").Value = -11084.667
# Row 84
$ws.Range("H84").Value = 2590.8235
$ws.Range("I84").Value = 2185.7144
$ws.Range("J84").Value = 4481.3335
$ws.Range("K84").Value = 21857.144
$ws.Range("L84").Value = 44813.335
$ws.Range("M84").Value = -16553.144
$ws.Range("N84").Value = -55421.335
# Row 96
$ws.Range("H96").Value = 3016.375
$ws.Range("J96").Value = 3221
$ws.Range("L96").Value = 3221
$ws.Range("N96").Value = -5967
